$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Api")

$ws.Range("C2").Value = "127"
$ws.Range("B2").Value = "285"
$ws.Range("C3").Value = "127"
$ws.Range("B3").Value = "286"
$ws.Range("C4").Value = "127"
$ws.Range("B4").Value = "287"

$ws.Range("B4").Select() | Out-Null
